$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 2088
$wsExhibit.Range("F5").Value = 1261
$wsExhibit.Range("F6").Value = 362

# Sheet "全部类型" updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 2088
$wsAll.Range("F7").Value = 1261
$wsAll.Range("F8").Value = 362
